$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# New log rows (MA cross with ML prediction) appended to the RSI log sheet.
$rows = @(
    @{ Row = 9;  Date = "2024-09-06"; Open = 76796000; Rsi = 30.54235498227352; Signal = 0;  Position = 0; Sell = $null;    High = $null },
    @{ Row = 10; Date = "2024-09-06"; Open = 76796000; Rsi = 30.54235498227352; Signal = 0;  Position = 0; Sell = $null;    High = $null },
    @{ Row = 11; Date = "2024-09-06"; Open = 76796000; Rsi = 30.54235498227352; Signal = 0;  Position = 0; Sell = $null;    High = $null },
    @{ Row = 12; Date = "2024-09-11"; Open = 77938000; Rsi = 31.95390905711808; Signal = 0;  Position = 1; Sell = $null;    High = 77938000 },
    @{ Row = 13; Date = "2024-09-11"; Open = 77938000; Rsi = 31.95390905711808; Signal = 0;  Position = 1; Sell = $null;    High = 77938000 },
    @{ Row = 14; Date = "2024-09-29"; Open = 86259000; Rsi = 69.91787829833063; Signal = -1; Position = 0; Sell = 86702330; High = 87001000 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Force column A to stay plain text (the date is stored as a string,
    # not an Excel date value) without leaving a lingering custom style.
    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.Date
    $cellA.Style = "Normal"

    $ws.Cells.Item($rowNum, 2).Value = $r.Open
    $ws.Cells.Item($rowNum, 3).Value = $r.Rsi
    $ws.Cells.Item($rowNum, 4).Value = $r.Signal
    $ws.Cells.Item($rowNum, 5).Value = $r.Position

    if ($r.Sell -ne $null) {
        $ws.Cells.Item($rowNum, 7).Value = $r.Sell
    }
    if ($r.High -ne $null) {
        $ws.Cells.Item($rowNum, 8).Value = $r.High
    }
}
